$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 1037.265132737054
$ws.Range("E2").Value = 28926.05393052954
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 16171.06685703679
$ws.Range("L2").Value = 48492.22142001599
$ws.Range("M2").Value = 10595.37713982
$ws.Range("N2").Value = 7071.74531360843
$ws.Range("O2").Value = 6993.890772562212

$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 4157.588990853394
$ws.Range("E2").Value = 45991.90904307188
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 37079.12819938764
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 17449.04999683176
$ws.Range("N2").Value = 9025.389658435586
$ws.Range("O2").Value = 9724.983840854013

$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13034.96639789043
$ws.Range("O2").Value = 12860.89728144265

$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13152.52576651623
$ws.Range("O2").Value = 12860.89728144265

$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13601.74312066917
$ws.Range("O2").Value = 14937.85618588151

$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13601.74312066917
$ws.Range("O2").Value = 14937.85618588151
